$wb = $excel.ActiveWorkbook

# Updated "want to go" counts (column F) for the "展览" and "全部类型" sheets.
$updates = @{
    2  = 815
    5  = 43
    6  = 12252
    9  = 489
    10 = 430
    12 = 890
    13 = 13588
    14 = 13722
    19 = 1020
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# Row 22's F value differs between the two sheets.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F22").Value = 2917

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F22").Value = 2918
